$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row right after the last used row (row 94 -> row 95)
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = 1.948
$ws.Cells.Item($newRow, 2).Value = 0.084
$ws.Cells.Item($newRow, 3).Value = 1.738
